$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.078.17"
$ws.Range("E2").Value = "  -4.43%  "

$ws.Range("D3").Value = "3.288.04"
$ws.Range("E3").Value = "  -4.80%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'554.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.52%  "

$ws.Range("D6").Value = "'142.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.81%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "3.285.64"
$ws.Range("E8").Value = "  -4.85%  "

$ws.Range("D9").Value = "'0.478"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "

$ws.Range("D10").Value = "'7.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.49%  "

$ws.Range("D11").Value = "'0.118"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.97%  "

$ws.Range("D12").Value = "'0.405"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.16%  "

$ws.Range("D13").Value = "3.854.54"
$ws.Range("E13").Value = "  -4.56%  "

$ws.Range("E14").Value = "  +0.33%  "

$ws.Range("D15").Value = "'27.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.49%  "

$ws.Range("D16").Value = "3.290.64"
$ws.Range("E16").Value = "  -4.34%  "

$ws.Range("D17").Value = "'0.0000165"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.93%  "

$ws.Range("D18").Value = "60.130.19"
$ws.Range("E18").Value = "  -4.31%  "

$ws.Range("D19").Value = "'6.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.65%  "

$ws.Range("D20").Value = "'14.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.98%  "

$ws.Range("D21").Value = "'8.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.87%  "

$ws.Range("D22").Value = "'371.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.23%  "

$ws.Range("D23").Value = "'73.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.03%  "

$ws.Range("D24").Value = "'0.542"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.75%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "3.447.04"
$ws.Range("E26").Value = "  -3.90%  "

$ws.Range("D27").Value = "'0.0000102"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -11.53%  "

$ws.Range("E28").Value = "  -5.77%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "'7.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.16%  "

$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'7.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.69%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.96%  "

$ws.Range("D34").Value = "'22.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.78%  "

$ws.Range("D35").Value = "'1.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.37%  "

$ws.Range("D36").Value = "'5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.10%  "

$ws.Range("D37").Value = "'166.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").Value = "'6.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.20%  "

$ws.Range("E39").Value = "  -8.53%  "

$ws.Range("D40").Value = "'26.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -16.46%  "

$ws.Range("D41").Value = "3.324.57"
$ws.Range("E41").Value = "  -4.58%  "

$ws.Range("D42").Value = "'0.0732"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.69%  "

$ws.Range("D43").Value = "'41.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.39%  "

$ws.Range("D44").Value = "'0.746"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.93%  "

$ws.Range("D45").Value = "'4.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.56%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.00%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.76%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.347.07"
$ws.Range("E49").Value = "  -8.46%  "

$ws.Range("D50").Value = "'6.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.11%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0255"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.32%  "
